$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap "Periodo Mora" (column E) and "Valor Mora" (column F) values between
# rows 16 and 17 so that period 2403 / value 52000 comes first (row 16)
# and period 2402 / value 38133 comes second (row 17).
$ws.Range("E16").Value = "2403"
$ws.Range("F16").Value = 52000

$ws.Range("E17").Value = "2402"
$ws.Range("F17").Value = 38133
